# Hortaliza, Terminal La Palmera de La Serena - Pepino dulce
# Insert a new weekly block (3 rows: Primera/Segunda/Tercera, fecha 2022-02-04 / serial 44596)
# right before the existing row 189 block, pushing all subsequent rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 189 (rows 189.. shift down to 192..)
$ws.Rows("189:191").Insert()

# Row 189 - Primera
$ws.Range("A189").Value = 8
$ws.Range("B189").Value = "Terminal La Palmera de La Serena"
$ws.Range("C189").Value = "Coquimbo"
$ws.Range("D189").Value = 44596
$ws.Range("E189").Value = 4
$ws.Range("F189").Value = 100112043
$ws.Range("G189").Value = "Pepino dulce"
$ws.Range("H189").Value = "Cultivar IV Región"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 440
$ws.Range("K189").Value = 14500
$ws.Range("L189").Value = 15000
$ws.Range("M189").Value = 14750
$ws.Range("N189").Value = "$/bandeja 18 kilos"
$ws.Range("O189").Value = "Provincia de Limarí"
$ws.Range("P189").Value = 819
$ws.Range("Q189").Value = 18
$ws.Range("R189").Value = "Hortaliza"

# Row 190 - Segunda
$ws.Range("A190").Value = 8
$ws.Range("B190").Value = "Terminal La Palmera de La Serena"
$ws.Range("C190").Value = "Coquimbo"
$ws.Range("D190").Value = 44596
$ws.Range("E190").Value = 4
$ws.Range("F190").Value = 100112043
$ws.Range("G190").Value = "Pepino dulce"
$ws.Range("H190").Value = "Cultivar IV Región"
$ws.Range("I190").Value = "Segunda"
$ws.Range("J190").Value = 360
$ws.Range("K190").Value = 12500
$ws.Range("L190").Value = 13000
$ws.Range("M190").Value = 12750
$ws.Range("N190").Value = "$/bandeja 18 kilos"
$ws.Range("O190").Value = "Provincia de Limarí"
$ws.Range("P190").Value = 708
$ws.Range("Q190").Value = 18
$ws.Range("R190").Value = "Hortaliza"

# Row 191 - Tercera
$ws.Range("A191").Value = 8
$ws.Range("B191").Value = "Terminal La Palmera de La Serena"
$ws.Range("C191").Value = "Coquimbo"
$ws.Range("D191").Value = 44596
$ws.Range("E191").Value = 4
$ws.Range("F191").Value = 100112043
$ws.Range("G191").Value = "Pepino dulce"
$ws.Range("H191").Value = "Cultivar IV Región"
$ws.Range("I191").Value = "Tercera"
$ws.Range("J191").Value = 300
$ws.Range("K191").Value = 10500
$ws.Range("L191").Value = 11000
$ws.Range("M191").Value = 10750
$ws.Range("N191").Value = "$/bandeja 18 kilos"
$ws.Range("O191").Value = "Provincia de Limarí"
$ws.Range("P191").Value = 597
$ws.Range("Q191").Value = 18
$ws.Range("R191").Value = "Hortaliza"
